# Scheduled runner update: refresh leve profit calculations (columns H:N) across all sheets
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 75.5
$ws.Range("I2").Value = 75.5
$ws.Range("K2").Value = 75.5
$ws.Range("M2").Value = 37.5
$ws.Range("H12").Value = 451
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 402
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 402
$ws.Range("M12").Value = -330
$ws.Range("N12").Value = -742
$ws.Range("H15").Value = 926.3582
$ws.Range("I15").Value = 926.3582
$ws.Range("K15").Value = 2779.0746
$ws.Range("M15").Value = -2610.0746
$ws.Range("H32").Value = 37041224
$ws.Range("I32").Value = 111113710
$ws.Range("J32").Value = 4978.1665
$ws.Range("K32").Value = 111113710
$ws.Range("L32").Value = 4978.1665
$ws.Range("M32").Value = -111113384
$ws.Range("N32").Value = -5630.1665
$ws.Range("H40").Value = 11959.1
$ws.Range("I40").Value = 18526.834
$ws.Range("J40").Value = 2107.5
$ws.Range("K40").Value = 18526.834
$ws.Range("L40").Value = 2107.5
$ws.Range("M40").Value = -18351.834
$ws.Range("N40").Value = -2457.5
$ws.Range("H98").Value = 37289.41
$ws.Range("I98").Value = 1695.7142
$ws.Range("K98").Value = 1695.7142
$ws.Range("M98").Value = -197.7141999999999
$ws.Range("H116").Value = 2528.743
$ws.Range("I116").Value = 1843.1818
$ws.Range("J116").Value = 3688.923
$ws.Range("K116").Value = 1843.1818
$ws.Range("L116").Value = 3688.923
$ws.Range("M116").Value = 1598.8182
$ws.Range("N116").Value = -10572.923
$ws.Range("H122").Value = 37289.41
$ws.Range("I122").Value = 1695.7142
$ws.Range("K122").Value = 5087.142599999999
$ws.Range("M122").Value = -2637.142599999999
$ws.Range("H137").Value = 5388.4736
$ws.Range("I137").Value = 4398.4614
$ws.Range("J137").Value = 7533.5
$ws.Range("K137").Value = 13195.3842
$ws.Range("L137").Value = 22600.5
$ws.Range("M137").Value = -10645.3842
$ws.Range("N137").Value = -27700.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13095.236
$ws.Range("I32").Value = 11963.878
$ws.Range("K32").Value = 11963.878
$ws.Range("M32").Value = -11676.878
$ws.Range("H61").Value = 1912.1482
$ws.Range("I61").Value = 1219.3889
$ws.Range("J61").Value = 3297.6667
$ws.Range("K61").Value = 1219.3889
$ws.Range("L61").Value = 3297.6667
$ws.Range("M61").Value = -1007.3889
$ws.Range("N61").Value = -3721.6667
$ws.Range("H74").Value = 2095.1892
$ws.Range("I74").Value = 1678.3793
$ws.Range("J74").Value = 3606.125
$ws.Range("K74").Value = 1678.3793
$ws.Range("L74").Value = 3606.125
$ws.Range("M74").Value = -804.3793000000001
$ws.Range("N74").Value = -5354.125
$ws.Range("H77").Value = 2095.1892
$ws.Range("I77").Value = 1678.3793
$ws.Range("J77").Value = 3606.125
$ws.Range("K77").Value = 8391.896500000001
$ws.Range("L77").Value = 18030.625
$ws.Range("M77").Value = -4023.896500000001
$ws.Range("N77").Value = -26766.625
$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -24992
$ws.Range("H136").Value = 1912.1482
$ws.Range("I136").Value = 1219.3889
$ws.Range("J136").Value = 3297.6667
$ws.Range("K136").Value = 3658.1667
$ws.Range("L136").Value = 9893.000100000001
$ws.Range("M136").Value = -1108.1667
$ws.Range("N136").Value = -14993.0001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 18000
$ws.Range("J40").Value = 18000
$ws.Range("L40").Value = 18000
$ws.Range("N40").Value = -18530
$ws.Range("H96").Value = 14589.5
$ws.Range("I96").Value = 4874.4
$ws.Range("J96").Value = 19986.777
$ws.Range("K96").Value = 4874.4
$ws.Range("L96").Value = 19986.777
$ws.Range("M96").Value = -2128.4
$ws.Range("N96").Value = -25478.777
$ws.Range("H97").Value = 4192.5
$ws.Range("I97").Value = 4192.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 4192.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -3201.5
$ws.Range("H134").Value = 2504.0637
$ws.Range("I134").Value = 2106.6875
$ws.Range("J134").Value = 3351.8
$ws.Range("K134").Value = 6320.0625
$ws.Range("L134").Value = 10055.4
$ws.Range("M134").Value = -3785.0625
$ws.Range("N134").Value = -15125.4
$ws.Range("N97").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3835022.8
$ws.Range("I31").Value = 1817.4706
$ws.Range("J31").Value = 6294060
$ws.Range("K31").Value = 1817.4706
$ws.Range("L31").Value = 6294060
$ws.Range("M31").Value = -1522.4706
$ws.Range("N31").Value = -6294650
$ws.Range("H34").Value = 3835022.8
$ws.Range("I34").Value = 1817.4706
$ws.Range("J34").Value = 6294060
$ws.Range("K34").Value = 1817.4706
$ws.Range("L34").Value = 6294060
$ws.Range("M34").Value = -1615.4706
$ws.Range("N34").Value = -6294464
$ws.Range("H58").Value = 2222.0312
$ws.Range("I58").Value = 1320.7778
$ws.Range("J58").Value = 3380.7856
$ws.Range("K58").Value = 1320.7778
$ws.Range("L58").Value = 3380.7856
$ws.Range("M58").Value = -1117.7778
$ws.Range("N58").Value = -3786.7856
$ws.Range("H93").Value = 4202.3335
$ws.Range("I93").Value = 4202.3335
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 4202.3335
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2330.3335
$ws.Range("H107").Value = 576.6111
$ws.Range("I107").Value = 379.52173
$ws.Range("J107").Value = 925.3077
$ws.Range("K107").Value = 379.52173
$ws.Range("L107").Value = 925.3077
$ws.Range("M107").Value = 1540.47827
$ws.Range("N107").Value = -4765.3077
$ws.Range("H132").Value = 563642.0600000001
$ws.Range("I132").Value = 2621.8667
$ws.Range("J132").Value = 1405172.2
$ws.Range("K132").Value = 7865.6001
$ws.Range("L132").Value = 4215516.6
$ws.Range("M132").Value = -5335.6001
$ws.Range("N132").Value = -4220576.6
$ws.Range("H134").Value = 286840.1
$ws.Range("I134").Value = 960.975
$ws.Range("J134").Value = 1557413.9
$ws.Range("K134").Value = 2882.925
$ws.Range("L134").Value = 4672241.699999999
$ws.Range("M134").Value = -347.9250000000002
$ws.Range("N134").Value = -4677311.699999999
$ws.Range("H136").Value = 2222.0312
$ws.Range("I136").Value = 1320.7778
$ws.Range("J136").Value = 3380.7856
$ws.Range("K136").Value = 3962.3334
$ws.Range("L136").Value = 10142.3568
$ws.Range("M136").Value = -1412.3334
$ws.Range("N136").Value = -15242.3568
$ws.Range("N93").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2093.818
$ws.Range("J34").Value = 2691.5
$ws.Range("L34").Value = 8074.5
$ws.Range("N34").Value = -8242.5
$ws.Range("H68").Value = 1211.7261
$ws.Range("I68").Value = 929.44446
$ws.Range("J68").Value = 1304.1091
$ws.Range("K68").Value = 2788.33338
$ws.Range("L68").Value = 3912.3273
$ws.Range("M68").Value = -1977.33338
$ws.Range("N68").Value = -5534.3273
$ws.Range("H71").Value = 1211.7261
$ws.Range("I71").Value = 929.44446
$ws.Range("J71").Value = 1304.1091
$ws.Range("K71").Value = 8365.00014
$ws.Range("L71").Value = 11736.9819
$ws.Range("M71").Value = -4309.00014
$ws.Range("N71").Value = -19848.9819
$ws.Range("H107").Value = 688.04
$ws.Range("I107").Value = 330.43478
$ws.Range("J107").Value = 794.8570999999999
$ws.Range("K107").Value = 991.3043399999999
$ws.Range("L107").Value = 2384.5713
$ws.Range("M107").Value = 928.6956600000001
$ws.Range("N107").Value = -6224.5713
$ws.Range("H131").Value = 52809.93
$ws.Range("J131").Value = 33286.188
$ws.Range("L131").Value = 99858.56400000001
$ws.Range("N131").Value = -109938.564

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1705
$ws.Range("I122").Value = 1433
$ws.Range("J122").Value = 2249
$ws.Range("K122").Value = 4299
$ws.Range("L122").Value = 6747
$ws.Range("M122").Value = -1849
$ws.Range("N122").Value = -11647
$ws.Range("H132").Value = 3052.08
$ws.Range("I132").Value = 2041.8334
$ws.Range("J132").Value = 3984.6155
$ws.Range("K132").Value = 6125.5002
$ws.Range("L132").Value = 11953.8465
$ws.Range("M132").Value = -3595.5002
$ws.Range("N132").Value = -17013.8465

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3607.5
$ws.Range("I46").Value = 1127.8334
$ws.Range("J46").Value = 6397.125
$ws.Range("K46").Value = 1127.8334
$ws.Range("L46").Value = 6397.125
$ws.Range("M46").Value = -939.8334
$ws.Range("N46").Value = -6773.125
$ws.Range("H55").Value = 476.84
$ws.Range("I55").Value = 392.30768
$ws.Range("J55").Value = 568.4167
$ws.Range("K55").Value = 392.30768
$ws.Range("L55").Value = 568.4167
$ws.Range("M55").Value = -219.30768
$ws.Range("N55").Value = -914.4167
$ws.Range("H100").Value = 4312
$ws.Range("I100").Value = 3766.6667
$ws.Range("J100").Value = 4609.4546
$ws.Range("K100").Value = 3766.6667
$ws.Range("L100").Value = 4609.4546
$ws.Range("M100").Value = -3225.6667
$ws.Range("N100").Value = -5691.4546
$ws.Range("H132").Value = 5905.4
$ws.Range("I132").Value = 6074.625
$ws.Range("J132").Value = 5712
$ws.Range("K132").Value = 18223.875
$ws.Range("L132").Value = 17136
$ws.Range("M132").Value = -15693.875
$ws.Range("N132").Value = -22196
$ws.Range("H136").Value = 2579.4
$ws.Range("I136").Value = 1931.3846
$ws.Range("K136").Value = 5794.1538
$ws.Range("M136").Value = -3244.1538

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2298.3635
$ws.Range("I132").Value = 2158.75
$ws.Range("J132").Value = 2378.1428
$ws.Range("K132").Value = 6476.25
$ws.Range("L132").Value = 7134.428400000001
$ws.Range("M132").Value = -3946.25
$ws.Range("N132").Value = -12194.4284
$ws.Range("H136").Value = 313607.1
$ws.Range("I136").Value = 417647.03
$ws.Range("J136").Value = 1487.25
$ws.Range("K136").Value = 1252941.09
$ws.Range("L136").Value = 4461.75
$ws.Range("M136").Value = -1250391.09
$ws.Range("N136").Value = -9561.75
